$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.470616203057875
$ws.Range("H2").Value = 0.775862068965517
$ws.Range("J2").Value = 1.525
$ws.Range("K2").Value = 0.033394541366277
$ws.Range("L2").Value = -0.260079093151912
$ws.Range("M2").Value = 0.234310225530207
$ws.Range("N2").Value = 2.18980599123128
$ws.Range("P2").Value = 'As likely as not improving'
$ws.Range("F3").Value = 0.5
$ws.Range("J3").Value = 11
$ws.Range("K3").Value = 0.0012722967821651
$ws.Range("L3").Value = -0.0597473741615733
$ws.Range("M3").Value = 0.0591652680288853
$ws.Range("N3").Value = 0.0115663343833198
$ws.Range("P3").Value = 'As likely as not increasing'
$ws.Range("E4").Value = 'ok'
$ws.Range("F4").Value = 0.0608830012402253
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.214285714285714
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0.0002502583764919
$ws.Range("M4").Value = 0.0005929383116882999
$ws.Range("N4").Value = 3.57511966417042
$ws.Range("P4").Value = 'Very unlikely improving'
$ws.Range("D5").Value = $true
$ws.Range("F5").Value = 0.999842716067363
$ws.Range("G5").Value = 0.0892857142857143
$ws.Range("H5").Value = 0.410714285714286
$ws.Range("J5").Value = 21
$ws.Range("K5").Value = -4.71290322580645
$ws.Range("L5").Value = -9.44185928167866
$ws.Range("M5").Value = -3.0103021978022
$ws.Range("N5").Value = -22.4423963133641
$ws.Range("P5").Value = 'Virtually certain improving'
$ws.Range("G6").Value = 0.981132075471698
$ws.Range("H6").Value = 0.0377358490566038
$ws.Range("I6").Value = 1
$ws.Range("F7").Value = 0.988347249152
$ws.Range("G7").Value = 0.892857142857143
$ws.Range("H7").Value = 0.0535714285714286
$ws.Range("P7").Value = 'Extremely likely improving'
$ws.Range("F8").Value = 0.980811510567775
$ws.Range("H8").Value = 0.714285714285714
$ws.Range("J8").Value = 0.0315
$ws.Range("K8").Value = -0.0036792582417582
$ws.Range("L8").Value = -0.0064740645563955
$ws.Range("M8").Value = -0.002124837206916
$ws.Range("N8").Value = -11.6801848944706
$ws.Range("P8").Value = 'Extremely likely improving'
$ws.Range("E9").Value = 'WARNING: Sen slope based on tied non-censored values'
$ws.Range("F9").Value = 0.476834566810564
$ws.Range("H9").Value = 0.781818181818182
$ws.Range("J9").Value = 7.55
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = -0.0468159264339306
$ws.Range("M9").Value = 0.0356879931658108
$ws.Range("N9").Value = 0
$ws.Range("P9").Value = 'As likely as not increasing'
$ws.Range("F10").Value = 0.998169198612947
$ws.Range("J10").Value = 0.03875
$ws.Range("K10").Value = -0.0040137362637362
$ws.Range("L10").Value = -0.007454081632653
$ws.Range("M10").Value = -0.002592682957231
$ws.Range("N10").Value = -10.3580290677065
$ws.Range("P10").Value = 'Virtually certain improving'
$ws.Range("E11").Value = 'ok'
$ws.Range("F11").Value = 0.992346212865899
$ws.Range("G11").Value = 0.0357142857142857
$ws.Range("H11").Value = 0.25
$ws.Range("K11").Value = -0.0052178571428571
$ws.Range("L11").Value = -0.0089197646631765
$ws.Range("N11").Value = -7.45408163265306
$ws.Range("P11").Value = 'Virtually certain improving'
$ws.Range("F12").Value = 0.0827352242598377
$ws.Range("K12").Value = 0.0002473140277069
$ws.Range("M12").Value = 0.0006604882459312
$ws.Range("N12").Value = 2.74793364118791
$ws.Range("P12").Value = 'Very unlikely improving'
$ws.Range("F13").Value = 0.013743168055755
$ws.Range("J13").Value = 0.63
$ws.Range("K13").Value = -0.0123488157040883
$ws.Range("M13").Value = -0.009859001579305899
$ws.Range("N13").Value = -1.96012947683941
$ws.Range("F14").Value = 0.04320536648685
$ws.Range("K14").Value = -4.88784233162507
$ws.Range("M14").Value = -2.64786529564445
$ws.Range("N14").Value = -3.87462729419348
$ws.Range("P14").Value = 'Extremely unlikely improving'
$ws.Range("F15").Value = 0.04320536648685
$ws.Range("K15").Value = -0.221433597781218
$ws.Range("M15").Value = -0.191983625505063
$ws.Range("N15").Value = -3.00452642851042
$ws.Range("P15").Value = 'Extremely unlikely improving'
